$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 233

# C2/D2/E2 need to stay plain text ("323", "2025-03-31", "2025-03-16") rather
# than being auto-coerced to a number / date by the usual Excel type
# inference. Force text mode via a temporary "@" (text) number format, enter
# the values, then clear the formatting back off again so the cells end up
# with no explicit style (matching a plain, never-formatted text cell).
$ws.Range("C2:E2").NumberFormat = "@"
$ws.Range("C2").Value = "323"
$ws.Range("D2").Value = "2025-03-31"
$ws.Range("E2").Value = "2025-03-16"
$ws.Range("C2:E2").ClearFormats()

$ws.Range("F2").Value = 1
